# Commit: "Tweak respiratory volume handling. Lots of documentation
# updates in preperation for release."
#
# The OOXML diff shows slide id="266" (the original first slide, a
# flowchart with "Preprocess" / "Process" / "PostProcess" boxes) removed
# from <p:sldIdLst>, with every other slide shifting up to take its
# place (old slide 2 -> new slide 1, old slide 3 -> new slide 2,
# old slide 4 -> new slide 3). No shape content on the surviving slides
# changes. So the edit is simply: delete the first slide of the deck.

$p = $ppt.ActivePresentation

# Slide 1 is the "Preprocess / Process / PostProcess" flowchart slide
# that disappears from the sldIdLst in the target revision.
$s = $p.Slides.Item(1)
$s.Delete()
